$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value is a plain numeric literal need to be forced
# to Text format first (otherwise Excel auto-converts "0.150" -> 0.15, etc.),
# then the cell style is reset back to Normal so no stray number format
# is left behind on the cell.

$ws.Range("D2").Value = '64.039.79'
$ws.Range("E2").Value = '  -1.49%  '
$ws.Range("D3").Value = '3.150.83'
$ws.Range("E3").Value = '  -0.93%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.11%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '3.146.58'
$ws.Range("E8").Value = '  -0.91%  '
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.150'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.39'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000255'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.98'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.67%  '
$ws.Range("D15").Value = '3.663.76'
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("E16").Value = '  +2.53%  '
$ws.Range("D17").Value = '64.007.97'
$ws.Range("E17").Value = '  -1.50%  '
$ws.Range("D18").Value = '3.146.74'
$ws.Range("E18").Value = '  -0.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.85'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '489.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.71'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.713'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.77'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.77%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.38%  '
$ws.Range("E28").Value = '  -5.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.99'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.08'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.79%  '
$ws.Range("E31").Value = '  +3.33%  '
$ws.Range("E32").Value = '  -5.52%  '
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("E34").Value = '  -2.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").Value = '0.0₃0754'
$ws.Range("E37").Value = '  -5.65%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.70'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.95'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0399'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '433.82'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.50%  '
$ws.Range("E42").Value = '  -0.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("D44").Value = '2.935.47'
$ws.Range("E44").Value = '  +2.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.261'
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = '  -6.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.96%  '
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.86'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.76%  '
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.97'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.15%  '
